$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns B:V) between row 55 and row 56 ---
# (Column A holds the running index and must stay put on each row.)
for ($col = 2; $col -le 22; $col++) {
  $v55 = $ws.Cells.Item(55, $col).Value2
  $v56 = $ws.Cells.Item(56, $col).Value2
  $ws.Cells.Item(55, $col).Value = $v56
  $ws.Cells.Item(56, $col).Value = $v55
}

# --- Append new row 78: Ararat-Armenia 2 x 0 Shirak Gyumri ---
$ws.Range("A78").Value = 77
$ws.Range("A77").Copy()
$ws.Range("A78").PasteSpecial(-4122)   # xlPasteFormats - match the index-column style

$ws.Range("B78").Value = "armenia"
$ws.Range("C78").Value = "premier-league"
$ws.Range("D78").Value = "2023-2024"

$ws.Range("E78").Value = 45241.625
$ws.Range("E77").Copy()
$ws.Range("E78").PasteSpecial(-4122)   # xlPasteFormats - match the date-column style

$ws.Range("F78").Value = "Ararat-Armenia"
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = "Shirak Gyumri"
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1.21
$ws.Range("K78").Value = "10/11/2023 03:13"
$ws.Range("L78").Value = 1.28
$ws.Range("M78").Value = "11/11/2023 09:33"
$ws.Range("N78").Value = 5.94
$ws.Range("O78").Value = "10/11/2023 03:13"
$ws.Range("P78").Value = 5.7
$ws.Range("Q78").Value = "11/11/2023 14:51"
$ws.Range("R78").Value = 9.83
$ws.Range("S78").Value = "10/11/2023 03:13"
$ws.Range("T78").Value = 10.09
$ws.Range("U78").Value = "11/11/2023 14:51"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/armenia/premier-league/ararat-armenia-shirak-gyumri/YisIzklU/"
